$wb = $excel.ActiveWorkbook
$new = $wb.Worksheets.Add()
$new.Name = "SheetNew"
$target = $wb.Worksheets.Item("SheetNew")
$target.Outline.SummaryRow = 1
$target.Outline.SummaryColumn = 1

$target.Range("A1").Value = "Модель"
$target.Range("B1").Value = "Тип коробу"
$target.Range("C1").Value = "Відкривання"
$target.Range("D1").Value = "Сторона"
$target.Range("E1").Value = "Полотно"
$target.Range("G1").Value = "Короб"
$target.Range("I1").Value = "Отвір"
$target.Range("K1").Value = "Оздоблення полотна"
$target.Range("L1").Value = "Алюм обв'язок"
$target.Range("M1").Value = "Колір фарбування профілю"
$target.Range("N1").Value = "Колір ущільнювача"
$target.Range("O1").Value = "Врізання"
$target.Range("R1").Value = "Колір фурн"
$target.Range("S1").Value = "к-сть"
$target.Range("T1").Value = "Ціна"
$target.Range("U1").Value = "Усього"
$target.Range("O2").Value = "Отвори"
$target.Range("Q2").Value = "Петлі"
$target.Range("E3").Value = "Ширина мм обличчя\тил"
$target.Range("F3").Value = "Висота мм обличчя\тил"
$target.Range("G3").Value = "Ширина мм"
$target.Range("H3").Value = "Висота мм"
$target.Range("I3").Value = "Ширина мм"
$target.Range("J3").Value = "Висота мм"
$target.Range("O3").Value = "Ручка"
$target.Range("P3").Value = "WC/PZ"
$target.Range("T3").Value = "грн"
$target.Range("U3").Value = "грн"
$target.Range("A4").Value = "Grezza PN ґрунт"
$target.Range("B4").Value = "SlimTS"
$target.Range("C4").Value = "Левое"
$target.Range("D4").Value = "лицьова"
$target.Range("E4").Value = "615"
$target.Range("F4").Value = "2100"
$target.Range("G4").Value = "689"
$target.Range("H4").Value = "2145"
$target.Range("I4").Value = "669"
$target.Range("J4").Value = "2135"
$target.Range("K4").Value = "Ґрунт"
$target.Range("L4").Value = "+"
$target.Range("M4").Value = "Чорний"
$target.Range("N4").Value = "Чорний"
$target.Range("O4").Value = "+"
$target.Range("P4").Value = "WC"
$target.Range("Q4").Value = "2"
$target.Range("R4").Value = "Чорний"
$target.Range("S4").Value = "5"
$target.Range("T4").Value = "15990"
$target.Range("U4").Value = "79950"
$target.Range("C5").Value = "Зовнішнє"
$target.Range("D5").Value = "тил"
$target.Range("E5").Value = "593"
$target.Range("F5").Value = "2089"
$target.Range("A9").Value = "Усього за дверні блоки: 79950.0 грн`n`t`tЗнижка:0.0%`n`t`tУсього, з урахуванням знижки: 79950.0 грн`n`t`tДоставка на склад (об'єкт) без вивантаження та занесення на поверх: 0.0 грн`n`t`tМонтаж:500.0грнЗаміри: 0.0 грн`n`t`tВсього за послуги: 2500.0 грн`n`t`tПідсумки без ПДВ: 82450.0 грн`n`t`tПередплата: 0.0 % `n`t`tПередплата: 0.0 грн `n`t`tЗалишок: 82450.0 грн"
$target.Range("H9").Value = "Місто: Замовник: `n`t`tДоставка замовлення: `n`t`tКонтакти: `n`t`tЕл.Адреса:"
$target.Range("O9").Value = "Висота ручки: 1000 мм. від низу полотна.`n`t`t`tЗазор від підлоги до полотна: 8 мм.`n`t`t`tБез вартості ручок.`n`t`t`tТермін виготовлення: 1-4 тижнів з моменту погодження та внесення передоплати.`n`t`t`tПримітка:"

$target.Range("N1:N3").Merge()
$target.Range("O1:Q1").Merge()
$target.Range("R4:R5").Merge()
$target.Range("T4:T5").Merge()
$target.Range("O4:O5").Merge()
$target.Range("A9:G18").Merge()
$target.Range("D1:D3").Merge()
$target.Range("Q4:Q5").Merge()
$target.Range("J4:J5").Merge()
$target.Range("G1:H2").Merge()
$target.Range("L4:L5").Merge()
$target.Range("I1:J2").Merge()
$target.Range("S1:S3").Merge()
$target.Range("T1:T2").Merge()
$target.Range("O9:U18").Merge()
$target.Range("A1:A3").Merge()
$target.Range("R1:R3").Merge()
$target.Range("H9:N18").Merge()
$target.Range("K1:K3").Merge()
$target.Range("A4:A5").Merge()
$target.Range("G4:G5").Merge()
$target.Range("I4:I5").Merge()
$target.Range("K4:K5").Merge()
$target.Range("S4:S5").Merge()
$target.Range("Q2:Q3").Merge()
$target.Range("U4:U5").Merge()
$target.Range("M1:M3").Merge()
$target.Range("U1:U2").Merge()
$target.Range("O2:P2").Merge()
$target.Range("B1:B3").Merge()
$target.Range("B4:B5").Merge()
$target.Range("C1:C3").Merge()
$target.Range("M4:M5").Merge()
$target.Range("H4:H5").Merge()
$target.Range("E1:F2").Merge()
$target.Range("N4:N5").Merge()
$target.Range("P4:P5").Merge()
$target.Range("L1:L3").Merge()

$old = $wb.Worksheets.Item("Sheet")
$old.Delete()
$wb.Worksheets.Item("SheetNew").Name = "Sheet"
